# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de), the two source-file rows (row 2 and
# row 3, corresponding to the two *.md files being localized) move from
# "Ready for handoff" to "Handed back: in sync with en-US". As part of
# recording the handback, the report now also fills in:
#   - Column E "Latest Target File"   -> same file/link as the source (col A)
#   - Column F "Latest Handback File" -> same file/link as the handoff xlf (col C)
#   - Column G "Latest Handback DateTime" -> the timestamp of the handback
#
# The status text itself ("Ready for handoff" -> "Handed back: in sync with
# en-US") is shared across the Overview sheet and both locale sheets via the
# shared string table, so updating it once on a locale sheet is reflected
# everywhere that text is used.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$sheetsInfo = @(
    @{
        Name = "zh-cn"
        HandbackTime = "2016-03-08 10:09:17"
        Rows = @(
            @{
                Row = 2
                MdDisplay = "7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.md"
                MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/c03d745db5bbdd77b2e394d21237938c99fc050e/e2e/7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.md"
                XlfDisplay = "7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.d3094665b3c23a473d3c8219d3415372c8b857e6.zh-cn.xlf"
                XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/606b3f7a710ac325e18688b8f12088070c35db91/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.d3094665b3c23a473d3c8219d3415372c8b857e6.zh-cn.xlf"
            },
            @{
                Row = 3
                MdDisplay = "9ceb1f0b-63eb-4439-a511-945c59e51e61.md"
                MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/c03d745db5bbdd77b2e394d21237938c99fc050e/e2e/9ceb1f0b-63eb-4439-a511-945c59e51e61.md"
                XlfDisplay = "9ceb1f0b-63eb-4439-a511-945c59e51e61.ff47824b987c0a7175e92c24c591a23ef598c876.zh-cn.xlf"
                XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/606b3f7a710ac325e18688b8f12088070c35db91/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/9ceb1f0b-63eb-4439-a511-945c59e51e61.ff47824b987c0a7175e92c24c591a23ef598c876.zh-cn.xlf"
            }
        )
    },
    @{
        Name = "de-de"
        HandbackTime = "2016-03-08 10:09:25"
        Rows = @(
            @{
                Row = 2
                MdDisplay = "7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.md"
                MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/c03d745db5bbdd77b2e394d21237938c99fc050e/e2e/7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.md"
                XlfDisplay = "7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.d3094665b3c23a473d3c8219d3415372c8b857e6.de-de.xlf"
                XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b25e96cb7aaa43ea94e1a0a7edd5cc58a3563a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/7ed6a9ed-d357-4dc9-a7ca-162ad1565ff3.d3094665b3c23a473d3c8219d3415372c8b857e6.de-de.xlf"
            },
            @{
                Row = 3
                MdDisplay = "9ceb1f0b-63eb-4439-a511-945c59e51e61.md"
                MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/c03d745db5bbdd77b2e394d21237938c99fc050e/e2e/9ceb1f0b-63eb-4439-a511-945c59e51e61.md"
                XlfDisplay = "9ceb1f0b-63eb-4439-a511-945c59e51e61.ff47824b987c0a7175e92c24c591a23ef598c876.de-de.xlf"
                XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b25e96cb7aaa43ea94e1a0a7edd5cc58a3563a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/9ceb1f0b-63eb-4439-a511-945c59e51e61.ff47824b987c0a7175e92c24c591a23ef598c876.de-de.xlf"
            }
        )
    }
)

foreach ($sheetInfo in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    foreach ($rowInfo in $sheetInfo.Rows) {
        $r = $rowInfo.Row

        # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Cells.Item($r, 2).Value = $newStatus

        # Latest Target File (E) - mirrors the source markdown file (col A)
        $eCell = $ws.Cells.Item($r, 5)
        $eCell.Value = $rowInfo.MdDisplay
        $ws.Hyperlinks.Add($eCell, $rowInfo.MdAddress, "", "", $rowInfo.MdDisplay) | Out-Null
        $eCell.Font.Underline = 2
        $eCell.Font.Color = 15570276

        # Latest Handback File (F) - mirrors the handoff xlf file (col C)
        $fCell = $ws.Cells.Item($r, 6)
        $fCell.Value = $rowInfo.XlfDisplay
        $ws.Hyperlinks.Add($fCell, $rowInfo.XlfAddress, "", "", $rowInfo.XlfDisplay) | Out-Null
        $fCell.Font.Underline = 2
        $fCell.Font.Color = 15570276

        # Latest Handback DateTime (G)
        $ws.Cells.Item($r, 7).Value = $sheetInfo.HandbackTime
    }
}

$wb.Save()
